$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for new rows 230-233 (dates 44304-44307, i.e. 2021-04-18 .. 2021-04-21)
$data = @(
    @(44304, 1, 12, 75.14088916718849),
    @(44305, 0, 7, 43.83218534752661),
    @(44306, 1, 8, 50.09392611145898),
    @(44307, 0, 6, 37.57044458359425)
)

$startRow = 230
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r - 1, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
}
